$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the two runs that make up the last paragraph ("Ceci n'est
#    pas sousligner ... {{UNDER_C_...}}.") into a single run/<w:t>.
#    Doing a Find/Replace across the run boundary (without changing the
#    visible text) makes Word rebuild that stretch of the paragraph as
#    one run, which is exactly the merge the diff shows.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "mais il est partie {{UNDER_C",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "mais il est partie {{UNDER_C",
    2
) | Out-Null

# ---------------------------------------------------------------------
# 2. Paragraph-format cleanup on the "Normal" style (if-tag / ODT
#    support pass): explicit left alignment, zeroed space before/after,
#    hyphenation suppressed, hanging punctuation turned off.
# ---------------------------------------------------------------------
$normal = $d.Styles("Normal")
$pf = $normal.ParagraphFormat
$pf.Hyphenation = $false
$pf.HangingPunctuation = $false
$pf.Alignment = 0
$pf.SpaceBefore = 0
$pf.SpaceAfter = 0
